$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New submission row appended to the FormData sheet (dimension grows A1:H5 -> A1:H6)
$ws.Range("A6").Value = "Travel"
$ws.Range("B6").Value = "1-2 Weeks"
$ws.Range("C6").Value = "rktindia2003@gmail.com"
$ws.Range("D6").Value = "Ravikant Tiwari"
# contactNumber / countryCode are textual (leading "+", long digit strings) in the
# source data, so force text entry the same way a real user typing an apostrophe-
# prefixed value into Excel would, keeping them from being coerced into numbers.
$ws.Range("E6").Value = "'8744883594"
$ws.Range("F6").Value = "'+40"
$ws.Range("G6").Value = "11/18/2024, 11:26:48 PM"
# Column H (timestamp) is left blank for this row, matching the source edit.
